$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-7, columns A-E: replace measurement values (stored as text, like the originals) ---
$newData = @{
    2 = @("4520.739832","3396.5556349999993","7364.308857999999","16783.801179","22977.184621999997")
    3 = @("8674.385119999997","6575.919615000003","14572.916676999996","33520.63456499999","46143.18019299999")
    4 = @("17505.421296999997","13122.930255999989","29212.589455","66916.72769099998","92091.07091100003")
    5 = @("34993.51118099997","26284.58743699997","58428.931455","134900.58596100003","184947.50486099996")
    6 = @("70649.45734499997","52244.13258100003","117184.53777700006","268124.09309700003","367320.8603580001")
    7 = @("155505.491163","116712.30815800004","246883.17314300002","577086.9314040001","795402.6730360001")
}

$cols = @("A","B","C","D","E")
foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$r"
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $vals[$i]
    }
}

# --- Row 17, column F: the sample count sequence now skips 16 and goes straight to 17 ---
$ws.Range("F17").Value = 17

# --- New rows 18-22: extra "n" values with no other columns populated ---
for ($r = 18; $r -le 22; $r++) {
    $ws.Range("F$r").Value = $r
}
